$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate formatting of the last existing row (row 25) into the new row 26
$ws.Range("A25:G25").Copy()
$ws.Range("A26:G26").PasteSpecial(-4122)

# Column F ("End date") in this sheet stores its dates as literal text
# (e.g. "2025-08-24" in row 25), not as a numeric date serial. Force the
# new F26 cell to be text (not an auto-parsed date) while keeping the
# same date-style number format as the rest of column F.
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "2025-08-25"
$ws.Range("F25").Copy()
$x = $ws.Range("F26").PasteSpecial(-4122)

# Fill in the rest of the new row's values
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Get approval"
$ws.Range("C26").Value = "Not Started"
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = "2025-07-23"
$ws.Range("G26").Value = "Thakarkuldip"

# Row 25's Assignees also changed (Aishwarrya VP -> Thakarkuldip)
$ws.Range("G25").Value = "Thakarkuldip"

# Grow the table to include the new row
$lo = $ws.ListObjects.Item(1)
$x = $lo.Resize($ws.Range("A1:H26"))

# Update the active selection to match the author's final cursor position
$x = $ws.Range("B31").Select()
